# Update the two "Report Date" values on the Queries sheet (Start Date / End Date)
# from 11-02-2020 / 13-02-2020 to 19-03-2020 / 27-03-2020, keeping their
# existing (quote-prefixed text) cell style, and leave the selection on the
# cell that was last edited.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

$ws.Range("D2").Value = "'19-03-2020 00:00:00"
$ws.Range("E2").Value = "'27-03-2020 00:00:00"

[void]$ws.Range("E2").Select()
